# Conceded Goal Maps Fixed | Vitor Hugo Radar done
#
# Adds a small "goals conceded vs. Premier League opponent" breakdown
# (Team name + total conceded, computed from a goal-by-goal formula) in
# columns B:C starting at row 14 of the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block (rows 14-19, columns B & C) ------------------------
$ws.Range("B14").Value = "City"
$ws.Range("C14").Formula = "=8+5+5+3+1+1+1+1"

$ws.Range("B15").Value = "United"
$ws.Range("C15").Formula = "=11+10+2+1+1+1"

$ws.Range("B16").Value = "Tottenham"
$ws.Range("C16").Formula = "=13+7+5+4+2+1+1+1+1"

$ws.Range("B17").Value = "Chelsea"
$ws.Range("C17").Formula = "=14+10+4+3+2+1+1"

$ws.Range("B18").Value = "Liverpool"
$ws.Range("C18").Formula = "=11+10+6+3+3+2+1+1+1"

$ws.Range("B19").Value = "WestHam"
$ws.Range("C19").Formula = "=26+18+5+5+3+2+2+2+2+1"

# --- Column B needs to widen to fit the new, longer team names ---------
$ws.Columns.Item(2).ColumnWidth = 9.8

# --- Match the selection left behind by the edit (B14 anchor, B14:C19 highlighted) ---
$ws.Range("B14:C19").Select()

Write-Host "done"
